$d = $word.ActiveDocument

$map = @{
    "57×39=2223" = "70×75=5250"
    "36×48=1728" = "41×59=2419"
    "13×65=845"  = "61×20=1220"
    "87×53=4611" = "69×90=6210"
    "90×17=1530" = "71×38=2698"
    "81×89=7209" = "95×58=5510"
    "42×69=2898" = "11×83=913"
    "71×97=6887" = "88×73=6424"
    "37×23=851"  = "46×36=1656"
    "93×84=7812" = "93×19=1767"
    "45×68=3060" = "78×21=1638"
    "20×64=1280" = "23×69=1587"
    "35×22=770"  = "51×60=3060"
    "88×66=5808" = "59×40=2360"
    "96×32=3072" = "88×97=8536"
    "15×79=1185" = "62×92=5704"
    "30×79=2370" = "85×71=6035"
    "58×90=5220" = "26×63=1638"
    "54×44=2376" = "51×71=3621"
    "55×57=3135" = "28×25=700"
    "25×79=1975" = "20×54=1080"
    "43×12=516"  = "91×90=8190"
    "39×16=624"  = "16×53=848"
    "28×86=2408" = "22×89=1958"
    "43×74=3182" = "93×88=8184"
}

foreach ($old in $map.Keys) {
    $new = $map[$old]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
